# Monads.docx edit: expand the "Reactive Streams" clause in the Zippers
# bullet with a parenthetical describing location observers/observables.

$d = $word.ActiveDocument

$old = "Zippers: Aggregation / Addressing: Locations / Contexts. Parsing. Monads. Augmentation (navigation / transforms) Reactive Streams."
$new = "Zippers: Aggregation / Addressing: Locations / Contexts. Parsing. Monads. Augmentation (navigation / transforms) Reactive Streams (location observers / observables: paths / kinds paths dataflow signatures)."

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Output ("replaced: " + $found)
